# Hindalco prices workbook update: add a new day's entry (10-11-2025) at the
# top of the table, shifting all existing rows down by one, which also
# duplicates the final (oldest) row to fill the row that rolls off the
# bottom of the historical window.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new blank row at row 2, pushing all data rows down one ---
# (this naturally extends the used range / last row from 152 to 153 and
# carries every shifted cell's value + formatting down with it)
$ws.Rows("2:2").Insert()

# --- 2. Restore formatting on the freshly inserted (blank) row 2, since a ---
# newly inserted row picks up the header row's bold font / number format
# instead of matching the rest of the (non-bold) data rows
$ws.Range("A2:C2").NumberFormat = "General"
$ws.Range("D2").NumberFormat = "0.000"
$ws.Range("E2:F2").NumberFormat = "General"
$ws.Range("A2:F2").HorizontalAlignment = -4108
$ws.Range("A2:F2").VerticalAlignment = -4108
$ws.Range("A2:F2").Font.Bold = $false
$ws.Range("A2:F2").Borders.LineStyle = -4142

# --- 3. Populate the new top row with the latest day's price data ---
# (same description/grade/price/circular as the prior top row; only the
# reported Date changes to the new day)
$ws.Cells.Item(2, 1).Value = "'10-11-2025"
$ws.Cells.Item(2, 2).Value = "2. P0610 (99.85% min) /P1020/ EC Grade Ingot & Sow 99.7% (min) / Cast Bar"
$ws.Cells.Item(2, 3).Value = "P1020"
$ws.Cells.Item(2, 4).Value = 288.25
$ws.Cells.Item(2, 5).Value = "'05.11.2025"
$ws.Cells.Item(2, 6).Value = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-05-november-2025.pdf"

# --- 4. Rebuild the Circular Link hyperlinks ---
# Row-shifting via Insert() does not renumber the worksheet's stored
# hyperlink anchors, so the old hyperlink list (anchored to F2..F91) would
# now point at the wrong rows. Clear them all out and recreate one
# hyperlink per row from the (now correctly shifted) text already sitting
# in column F, wherever that text is non-empty.
$ws.Hyperlinks.Delete()

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $link = $ws.Cells.Item($r, 6).Value2
    if ($link -ne $null -and $link -ne "") {
        $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $link, [Type]::Missing, [Type]::Missing, $link)
    }
}

# Adding a hyperlink auto-applies Excel's default blue/underlined "Hyperlink"
# visual style; the source data keeps plain (non-hyperlink-styled) text, so
# restore the original look across every data row/column.
$dataRange = $ws.Range("A2:F" + $lastRow)
$dataRange.Font.Underline = -4142
$dataRange.Font.ThemeColor = 1
$dataRange.Font.TintAndShade = 0

Write-Host "Update complete. Last row: $lastRow"
